$d = $word.ActiveDocument

# --- Edit 1: "...Actuarial Science Program Central Washington University" ->
#             "...Actuarial Science Program at Central Washington University"
$d.Content.Find.Execute(
    "Actuarial Science Program Central Washington University",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Actuarial Science Program at Central Washington University", 2) | Out-Null

# --- Edit 2: "frustration. Those with" -> "frustrations. Those with"
$d.Content.Find.Execute(
    "minimal frustration. Those with",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "minimal frustrations. Those with", 2) | Out-Null

# --- Edit 3: "one user type, basic user." -> "one user type, the basic user."
$d.Content.Find.Execute(
    "one user type, basic user.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "one user type, the basic user.", 2) | Out-Null

# --- Edit 4: "it will be increasingly" -> "it will become increasingly"
#     and reposition the _GoBack bookmark to sit right before "increasingly"
$fr = $d.Content
$fr.Find.Execute("be increasingly", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$becomePos = $fr.Start + 2            # right after "be"
$insRange = $d.Range($becomePos, $becomePos)
$insRange.InsertAfter("come")

# Find the new location right before "increasingly" for the bookmark
$fr2 = $d.Content
$fr2.Find.Execute("increasingly important", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $fr2.Start

# Remove the existing _GoBack bookmark (currently sitting alone in the last paragraph)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Merge the final two paragraphs: delete the paragraph mark that currently
# separates the "...events." paragraph from the (now empty) bookmark paragraph.
$lastParaIndex = $d.Paragraphs.Count
$secondLastParaIndex = $lastParaIndex - 1
$mergeRange = $d.Range($d.Paragraphs.Item($secondLastParaIndex).Range.End - 1, $d.Paragraphs.Item($secondLastParaIndex).Range.End)
$mergeRange.Delete()

# Re-add the _GoBack bookmark right before "increasingly"
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Add a brand new empty paragraph at the very end of the document.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()

Write-Host "Done."
Write-Host $d.Content.Text
